# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    3  = 296
    4  = 11160
    5  = 10389
    6  = 593
    8  = 735
    10 = 19
    13 = 9628
    14 = 2219
    16 = 2443
    19 = 87
    20 = 392
    21 = 10861
    22 = 10793
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
